# Demography.xlsx update
#  - rename sheet 2 ("1618 lake sims parkway oc" -> "test")
#  - refresh the "Census track" (column C) demographic figures on the
#    "451 clear blue way mcdono" sheet
#  - bump the Google Maps link zoom/params on that same sheet (B13)

$wb = $excel.ActiveWorkbook

# --- rename the second sheet -------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "test"

# --- update the Census track column on the third sheet ----------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C2").Value  = "29,314 (100%)"
$ws3.Range("C3").Value  = "2.92%(per year)"
$ws3.Range("C4").Value  = "1.92%(per year)"

# these three look like plain currency ("$nn,nnn") so Excel's automatic
# type detection would otherwise coerce them to Currency-formatted numbers;
# force the cell to Text first so the literal string is preserved, exactly
# like a user would do when typing a "numeric-looking" label.
$ws3.Range("C6").NumberFormat = "@"
$ws3.Range("C6").Value = "$63,043"

$ws3.Range("C7").NumberFormat = "@"
$ws3.Range("C7").Value = "$84,267"

$ws3.Range("C8").Value  = "11,561 (100%)"
$ws3.Range("C9").Value  = "6,674 (57.7%)"
$ws3.Range("C10").Value = "3,973 (34.4%)"
$ws3.Range("C11").Value = "913 ( 7.9%)"

$ws3.Range("C12").NumberFormat = "@"
$ws3.Range("C12").Value = "$172,046"

# --- update the Google Maps link (zoom level + query params changed) -------
$ws3.Range("B13").Value = "https://www.google.com/maps/place/451+Clear+Blue+Way,+McDonough,+GA+30253,+USA/@33.4500017,-84.1641469,17z/data=!3m1!4b1!4m5!3m4!1s0x88f45a46fdfa89cb:0x798192984c4bcdd6!8m2!3d33.4500017!4d-84.1619582"
